# 自动更新Excel文件 - 2026-02-07 23:21:41
# Decrement the "剩余" (remaining) value in column E by 1 for every data row,
# except row 36 which stays unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 99; $row++) {
    if ($row -eq 36) {
        continue
    }
    $cell = $ws.Cells.Item($row, 5)  # Column E
    $current = $cell.Value2
    $cell.Value = $current - 1
}
